$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.020.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.926.34"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.90"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.72"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.94"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.410.45"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.931.74"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.926.40"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.71"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "432.19"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.08"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.42"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.89"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.21"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.91"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.26"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.16%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.01"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.65"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.57%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0859"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.02"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.98"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.83%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.62"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "380.48"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.702.05"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0342"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.01"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.80"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.12%  "
